$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 88. This shifts the existing rows 88-98
# down to rows 92-102, and creates empty rows 88-91 ready for new data.
$ws.Range("A88:A91").EntireRow.Insert()

# Common constant values for this block of rows.
$company = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$codreg = 15
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad = "Sin especificar"
$unidad = "`$/bandeja 3 kilos"
$origen = "Región de Arica y Parinacota"
$kgUnidad = 3

# New week of data (fecha 45154) inserted at rows 88-91.
$newRows = @(
    @{ Row = 88; Fecha = 45154; Calidad = "Especial"; Volumen = 200; PMin = 7000; PMax = 8000; PProm = 7500; PKg = 2500 },
    @{ Row = 89; Fecha = 45154; Calidad = "Primera";  Volumen = 200; PMin = 6000; PMax = 7000; PProm = 6500; PKg = 2167 },
    @{ Row = 90; Fecha = 45154; Calidad = "Segunda";  Volumen = 200; PMin = 4000; PMax = 5000; PProm = 4500; PKg = 1500 },
    @{ Row = 91; Fecha = 45154; Calidad = "Tercera";  Volumen = 200; PMin = 2000; PMax = 3000; PProm = 2500; PKg = 833 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = $company
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
